# TP3 Modelos y simulaciones
# Replace the trajectory data (rows 3-41) in columns A:D with the
# results of the re-run simulation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    (1.767766952966369, 0.4877258050403206, 1.669766952966369, 0),
    (3.535533905932737, 0.9754516100806413, 3.241533905932738, 0),
    (5.303300858899106, 1.463177415120962, 4.715300858899107, 0.1),
    (7.071067811865475, 1.950903220161283, 6.091067811865476, 0.2),
    (8.838834764831843, 2.438629025201603, 7.368834764831845, 0.3),
    (10.60660171779821, 2.926354830241924, 8.548601717798213, 0.4),
    (12.37436867076458, 3.414080635282245, 9.630368670764582, 0.5),
    (14.14213562373095, 3.901806440322566, 10.61413562373095, 0.6000000000000001),
    (15.04488810900642, 3.948785359893499, 11.49990257669732, 0.7000000000000001),
    (15.082626126591, 3.555017393995045, 12.28766952966369, 0.8),
    (14.25534967648468, 2.720502542627205, 12.97743648263006, 0.9),
    (13.42807322637836, 1.885987691259364, 13.56920343559642, 1),
    (12.60079677627204, 1.051472839891523, 14.06297038856279, 1.1),
    (11.77352032616572, 0.2169579885236828, 14.45873734152916, 1.2),
    (10.9462438760594, -0.6175568628441579, 14.75650429449553, 1.3),
    (10.11896742595308, -1.452071714211999, 14.9562712474619, 1.4),
    (9.291690975846761, -2.286586565579839, 15.05803820042827, 1.5),
    (8.464414525740441, -3.12110141694768, 15.06180515339463, 1.6),
    (7.637138075634122, -3.95561626831552, 14.967572106361, 1.7),
    (6.809861625527803, -4.790131119683361, 14.77533905932737, 1.8),
    (5.982585175421485, -5.624645971051201, 14.48510601229374, 1.9),
    (5.155308725315166, -6.459160822419042, 14.09687296526011, 2),
    (4.328032275208847, -7.293675673786883, 13.61063991822647, 2.1),
    (3.500755825102528, -8.128190525154723, 13.02640687119284, 2.2),
    (2.67347937499621, -8.962705376522564, 12.34417382415921, 2.3),
    (1.846202924889891, -9.797220227890405, 11.56394077712558, 2.4),
    (1.018926474783572, -10.63173507925825, 10.68570773009194, 2.5),
    (0.1916500246772526, -11.46624993062609, 9.709474683058312, 2.6),
    (-0.6356264254290664, -12.30076478199393, 8.635241636024681, 2.7),
    (-1.462902875535385, -13.13527963336177, 7.463008588991048, 2.8),
    (-2.290179325641704, -13.96979448472961, 6.192775541957416, 2.9),
    (-3.117455775748024, -14.80430933609745, 4.824542494923783, 3),
    (-3.944732225854342, -15.63882418746529, 3.35830944789015, 3.1),
    (-4.772008675960661, -16.47333903883313, 1.794076400856517, 3.2),
    (-5.59928512606698, -17.30785389020097, 0.1318433538228847, 3.3),
    (-6.426561576173299, -18.14236874156882, -1.628389693210748, 3.4),
    (-7.253838026279618, -18.97688359293666, -3.486622740244381, 3.5),
    (-8.081114476385936, -19.8113984443045, -5.442855787278013, 3.6),
    (-8.908390926492256, -20.64591329567234, -7.497088834311645, 3.7)
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
